$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: date 2023-02-07 (serial 44964) + progress note
$ws.Range("A2").Value = 44964
$ws.Range("A2").NumberFormat = "mm-dd-yy"
$ws.Range("B2").Value = "spring basics,login page"

# Row 3: date 2023-02-08 (serial 44965), re-using A2's date style, + progress note
$ws.Range("A2").Copy($ws.Range("A3"))
$ws.Range("A3").Value = 44965
$ws.Range("B3").Value = "c++ exception handling, working on login page"

# Column A width (matches bestFit width in target)
$ws.Columns("A").ColumnWidth = 9.5

# Selection moves to C6 like in the target sheetView
$ws.Range("C6").Select() | Out-Null
